# Update "想去人数" (F column) counts that changed between the two data
# refreshes. The same values need to be updated on both the "展览" sheet
# and the "全部类型" sheet, since they mirror the same underlying data.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 364
    "F3"  = 4731
    "F5"  = 422
    "F7"  = 939
    "F9"  = 2038
    "F11" = 1219
    "F13" = 83
    "F14" = 34
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
